$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("registrationData")

$ws.Range("A2").Value = "Firstname"
$ws.Range("B2").Value = "Lastname"
$ws.Range("C2").Value = "test6@gmail.com"
$ws.Range("E2").Value = "hello123"

$ws.Range("A3").Value = "Firstname"
$ws.Range("B3").Value = "Lastname"
$ws.Range("C3").Value = "test7@gmail.com"
$ws.Range("E3").Value = "hello123"

$ws.Range("A4").Value = "Firstname"
$ws.Range("B4").Value = "Lastname"
$ws.Range("C4").Value = "test8@gmail.com"
$ws.Range("E4").Value = "hello123"

$ws.Range("A5").Value = "Firstname"
$ws.Range("B5").Value = "Lastname"
$ws.Range("C5").Value = "test9@gmail.com"
$ws.Range("E5").Value = "hello123"

$ws.Range("A6").Value = "Firstname"
$ws.Range("B6").Value = "Lastname"
$ws.Range("C6").Value = "test10@gmail.com"
$ws.Range("E6").Value = "hello123"

$ws.Range("E5").Select()
